# Card13: "إضافة حدث جديد في Card13" (add a new event row)
# - Row 25's previously-blank B:K cells get filled with the "nan" placeholder
#   text used throughout this sheet for missing values.
# - A brand new row 26 is appended, duplicating row 25's event metadata
#   (date / correction / event / serviced-by) while leaving its own B:K
#   measurement columns blank, matching the existing data pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card13")

# --- Fix up row 25: B25:K25 were blank, now hold the "nan" placeholder ---
$naCols = @("B","C","D","E","F","G","H","I","J","K")
foreach ($col in $naCols) {
    $ws.Range($col + "25").Value = "nan"
}

# --- Append new row 26 (copy of row 25's event, fresh blank measurements) ---

# A26 holds the card number "13" as TEXT (matches the text-typed column A
# used throughout the sheet), so force text interpretation.
$ws.Range("A26").Value = "'13"

# B26:K26 stay blank/empty, same as row 25 was before this edit. Touch each
# cell with a no-op formatting call so an (empty) cell entry exists for the
# new row without giving it any value.
foreach ($col in $naCols) {
    $ws.Range($col + "26").Font.Bold = $false
}

# L26: date text (kept as literal text, not converted to a date serial)
$ws.Range("L26").Value = "17/11/2025"

# M26 / N26: correction + event description (Arabic text)
$ws.Range("M26").Value = "تم تغير سير  دوبل700(محمد نعيم)"
$ws.Range("N26").Value = "قطع سير كويلر مسنن دبل 700"

# O26: serviced by
$ws.Range("O26").Value = "فني"
